# Update the closing slide's title from "Le Fin" to "Thank You".
$p = $ppt.ActivePresentation

# The edited title ("Le Fin" -> "Thank You") lives on the last slide
# (slide 24) of the deck, in the Title placeholder shape.
$s = $p.Slides.Item(24)
$shp = $s.Shapes.Item(1)

$shp.TextFrame.TextRange.Text = "Thank You"
